# Update CDA Logical model for ST.r2b
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Include from ObservationInter")

# Rename the "Include from ObservationInter" tab to "Include #0"
$ws2.Name = "Include #0"

# Metadata sheet: bump Version and Date
$ws1.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws1.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row right after "Contact" (row 10),
# pushing Description/Purpose/Copyright/Immutable down by one row.
$ws1.Rows.Item(11).Insert()

# Copy formatting from the row below (the old row 11, now row 12) so the
# new row matches the table's existing style.
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)

$ws1.Range("A11").Value = "Jurisdiction"
$ws1.Range("B11").Value = ""

Write-Output "done"
